$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
$ws.Activate()

# Insert a new row above the old footer/note row (currently row 79),
# pushing that footer row (and its formatting) down to row 80.
$ws.Rows.Item(79).Insert()

# Fill in the new data row 79 with the latest day's figures.
$ws.Cells.Item(79, 1).Value = 43934
$ws.Cells.Item(79, 2).Value = 1047
$ws.Cells.Item(79, 3).Value = 22005
$ws.Cells.Item(79, 4).Value = 195
$ws.Cells.Item(79, 5).Value = 5067

# Update the print area to include the new last row.
$wb.Names.Item(1).RefersTo = "=相談件数!`$A`$1:`$E`$84"

# Update the active selection to reflect the new bottom of the sheet.
$ws.Range("G77").Select() | Out-Null
